$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for rows 3-6 per repulled data
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 1
